# Updates NATMI LR-pair metrics (Wnt5a-Fzd8) to reflect recomputed TPM values.
# Only the numeric value cells change; labels/columns A-F (for rows 5-10) stay the same.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: updated TPM-derived NATMI metrics
$ws.Range("G2").Value = 0.02354566666666667
$ws.Range("H2").Value = 0.07063700000000001
$ws.Range("I2").Value = 0.002815555392485919
$ws.Range("J2").Value = 0.002815555392485918
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.235341333333333
$ws.Range("N2").Value = 9.706023999999999
$ws.Range("O2").Value = 0.2153734454473681
$ws.Range("P2").Value = 0.2153734454473681
$ws.Range("Q2").Value = 0.07617826858755555
$ws.Range("R2").Value = 0.685604417288
$ws.Range("S2").Value = 0.0006063958657276091
$ws.Range("T2").Value = 0.000606395865727609

# Row 3: updated TPM-derived NATMI metrics
$ws.Range("G3").Value = 0.02354566666666667
$ws.Range("H3").Value = 0.07063700000000001
$ws.Range("I3").Value = 0.002815555392485919
$ws.Range("J3").Value = 0.002815555392485918
$ws.Range("O3").Value = 0.4841904166376352
$ws.Range("P3").Value = 0.4841904166376352
$ws.Range("Q3").Value = 0.1712596811994445
$ws.Range("R3").Value = 1.541337130795
$ws.Range("S3").Value = 0.001363264938554098
$ws.Range("T3").Value = 0.001363264938554097

# Row 4: updated TPM-derived NATMI metrics
$ws.Range("G4").Value = 0.02354566666666667
$ws.Range("H4").Value = 0.07063700000000001
$ws.Range("I4").Value = 0.002815555392485919
$ws.Range("J4").Value = 0.002815555392485918
$ws.Range("O4").Value = 0.3004361379149967
$ws.Range("P4").Value = 0.3004361379149967
$ws.Range("Q4").Value = 0.1062652118507778
$ws.Range("R4").Value = 0.956386906657
$ws.Range("S4").Value = 0.0008458945882042123
$ws.Range("T4").Value = 0.0008458945882042121

# Row 5: updated TPM-derived NATMI metrics
$ws.Range("I5").Value = 0.9868456480383168
$ws.Range("J5").Value = 0.9868456480383166
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.235341333333333
$ws.Range("N5").Value = 9.706023999999999
$ws.Range("O5").Value = 0.2153734454473681
$ws.Range("P5").Value = 0.2153734454473681
$ws.Range("Q5").Value = 26.70030681383555
$ws.Range("R5").Value = 240.30276132452
$ws.Range("S5").Value = 0.212540347342753
$ws.Range("T5").Value = 0.212540347342753

# Row 6: updated TPM-derived NATMI metrics
$ws.Range("I6").Value = 0.9868456480383168
$ws.Range("J6").Value = 0.9868456480383166
$ws.Range("O6").Value = 0.4841904166376352
$ws.Range("P6").Value = 0.4841904166376352
$ws.Range("S6").Value = 0.4778212054807097
$ws.Range("T6").Value = 0.4778212054807097

# Row 7: updated TPM-derived NATMI metrics
$ws.Range("I7").Value = 0.9868456480383168
$ws.Range("J7").Value = 0.9868456480383166
$ws.Range("O7").Value = 0.3004361379149967
$ws.Range("P7").Value = 0.3004361379149967
$ws.Range("S7").Value = 0.2964840952148541
$ws.Range("T7").Value = 0.296484095214854

# Row 8: updated TPM-derived NATMI metrics
$ws.Range("I8").Value = 0.0103387965691973
$ws.Range("J8").Value = 0.0103387965691973
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.235341333333333
$ws.Range("N8").Value = 9.706023999999999
$ws.Range("O8").Value = 0.2153734454473681
$ws.Range("P8").Value = 0.2153734454473681
$ws.Range("Q8").Value = 0.279728690127111
$ws.Range("R8").Value = 2.517558211143999
$ws.Range("S8").Value = 0.002226702238887452
$ws.Range("T8").Value = 0.002226702238887452

# Row 9: updated TPM-derived NATMI metrics
$ws.Range("I9").Value = 0.0103387965691973
$ws.Range("J9").Value = 0.0103387965691973
$ws.Range("O9").Value = 0.4841904166376352
$ws.Range("P9").Value = 0.4841904166376352
$ws.Range("S9").Value = 0.005005946218371395
$ws.Range("T9").Value = 0.005005946218371395

# Row 10: updated TPM-derived NATMI metrics
$ws.Range("I10").Value = 0.0103387965691973
$ws.Range("J10").Value = 0.0103387965691973
$ws.Range("O10").Value = 0.3004361379149967
$ws.Range("P10").Value = 0.3004361379149967
$ws.Range("R10").Value = 3.511878933640999
$ws.Range("S10").Value = 0.003106148111938456
$ws.Range("T10").Value = 0.003106148111938456

